$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (want-to-go count) column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 1925
$wsExpo.Range("F5").Value = 876
$wsExpo.Range("F6").Value = 270

# Sheet "全部类型" (all types) - same underlying events, update column F too
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 1925
$wsAll.Range("F6").Value = 876
$wsAll.Range("F7").Value = 270
